$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1347
$wsExpo.Range("F6").Value = 10570
$wsExpo.Range("F13").Value = 12418

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1347
$wsAll.Range("F7").Value = 10570
$wsAll.Range("F14").Value = 12418
